$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("less_10_inc")
$ws.Range("G2").Value = 3
$ws.Range("J2").Value = 26
$ws.Range("J3").Value = 38
$ws.Range("U3").Value = 42
$ws.Range("J4").Value = 52
$ws.Range("U4").Value = 56
$ws.Range("B5").Value = 24
$ws.Range("G5").Value = 8
$ws.Range("J5").Value = 57
$ws.Range("L5").Value = 516
$ws.Range("M5").Value = 734
$ws.Range("R5").Value = 877
$ws.Range("S5").Value = 568
$ws.Range("T5").Value = 290
$ws.Range("U5").Value = 63
$ws.Range("B6").Value = 24
$ws.Range("G6").Value = 8
$ws.Range("J6").Value = 59
$ws.Range("L6").Value = 563
$ws.Range("M6").Value = 809
$ws.Range("R6").Value = 953
$ws.Range("S6").Value = 622
$ws.Range("T6").Value = 315
$ws.Range("U6").Value = 66
$ws.Range("B7").Value = 24
$ws.Range("G7").Value = 8
$ws.Range("J7").Value = 59
$ws.Range("L7").Value = 563
$ws.Range("M7").Value = 809
$ws.Range("R7").Value = 953
$ws.Range("S7").Value = 622
$ws.Range("T7").Value = 315
$ws.Range("U7").Value = 66

$ws = $wb.Worksheets.Item("less_50_inc")
$ws.Range("F2").Value = 7
$ws.Range("J2").Value = 47
$ws.Range("L2").Value = 259
$ws.Range("F3").Value = 8
$ws.Range("J3").Value = 63
$ws.Range("K3").Value = 191
$ws.Range("L3").Value = 452
$ws.Range("T3").Value = 274
$ws.Range("U3").Value = 78
$ws.Range("F4").Value = 9
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 238
$ws.Range("L4").Value = 566
$ws.Range("T4").Value = 342
$ws.Range("U4").Value = 102
$ws.Range("B5").Value = 34
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 13
$ws.Range("J5").Value = 85
$ws.Range("K5").Value = 261
$ws.Range("L5").Value = 666
$ws.Range("T5").Value = 396
$ws.Range("U5").Value = 110
$ws.Range("B6").Value = 37
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 14
$ws.Range("J6").Value = 89
$ws.Range("K6").Value = 280
$ws.Range("L6").Value = 718
$ws.Range("T6").Value = 426
$ws.Range("U6").Value = 115
$ws.Range("B7").Value = 37
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 14
$ws.Range("J7").Value = 89
$ws.Range("K7").Value = 280
$ws.Range("L7").Value = 718
$ws.Range("T7").Value = 426
$ws.Range("U7").Value = 115

$ws = $wb.Worksheets.Item("great_10_inc")
$ws.Range("J2").Value = 73
$ws.Range("J3").Value = 108
$ws.Range("J4").Value = 127
$ws.Range("B5").Value = 55
$ws.Range("G5").Value = 23
$ws.Range("J5").Value = 141
$ws.Range("L5").Value = 341
$ws.Range("M5").Value = 281
$ws.Range("R5").Value = 215
$ws.Range("S5").Value = 312
$ws.Range("T5").Value = 324
$ws.Range("B6").Value = 55
$ws.Range("G6").Value = 23
$ws.Range("J6").Value = 150
$ws.Range("L6").Value = 355
$ws.Range("M6").Value = 289
$ws.Range("R6").Value = 222
$ws.Range("S6").Value = 322
$ws.Range("T6").Value = 340
$ws.Range("B7").Value = 55
$ws.Range("G7").Value = 23
$ws.Range("J7").Value = 150
$ws.Range("L7").Value = 355
$ws.Range("M7").Value = 289
$ws.Range("R7").Value = 222
$ws.Range("S7").Value = 322
$ws.Range("T7").Value = 340

$ws = $wb.Worksheets.Item("great_50_inc")
$ws.Range("J2").Value = 59
$ws.Range("L2").Value = 32
$ws.Range("J3").Value = 91
$ws.Range("K3").Value = 141
$ws.Range("L3").Value = 71
$ws.Range("T3").Value = 139
$ws.Range("U3").Value = 93
$ws.Range("J4").Value = 110
$ws.Range("K4").Value = 176
$ws.Range("L4").Value = 100
$ws.Range("T4").Value = 172
$ws.Range("U4").Value = 114
$ws.Range("B5").Value = 51
$ws.Range("G5").Value = 23
$ws.Range("J5").Value = 123
$ws.Range("K5").Value = 202
$ws.Range("L5").Value = 124
$ws.Range("T5").Value = 200
$ws.Range("U5").Value = 127
$ws.Range("B6").Value = 52
$ws.Range("G6").Value = 23
$ws.Range("J6").Value = 133
$ws.Range("K6").Value = 216
$ws.Range("L6").Value = 129
$ws.Range("T6").Value = 208
$ws.Range("U6").Value = 135
$ws.Range("B7").Value = 52
$ws.Range("G7").Value = 23
$ws.Range("J7").Value = 133
$ws.Range("K7").Value = 216
$ws.Range("L7").Value = 129
$ws.Range("T7").Value = 208
$ws.Range("U7").Value = 135

$ws = $wb.Worksheets.Item("less_10_exc")
$ws.Range("G2").Value = 3
$ws.Range("J2").Value = 26
$ws.Range("U3").Value = 27
$ws.Range("B5").Value = 11
$ws.Range("G5").Value = 4
$ws.Range("L5").Value = 102
$ws.Range("M5").Value = 145
$ws.Range("R5").Value = 163
$ws.Range("S5").Value = 123
$ws.Range("T5").Value = 66

$ws = $wb.Worksheets.Item("less_50_exc")
$ws.Range("F2").Value = 7
$ws.Range("J2").Value = 47
$ws.Range("L2").Value = 259
$ws.Range("F3").Value = 6
$ws.Range("K3").Value = 110
$ws.Range("L3").Value = 233
$ws.Range("T3").Value = 152
$ws.Range("U3").Value = 43
$ws.Range("K4").Value = 80
$ws.Range("L4").Value = 145
$ws.Range("T4").Value = 103
$ws.Range("U4").Value = 45
$ws.Range("B5").Value = 16
$ws.Range("G5").Value = 8
$ws.Range("L5").Value = 119
$ws.Range("T5").Value = 87

$ws = $wb.Worksheets.Item("great_10_exc")
$ws.Range("J2").Value = 73
$ws.Range("B5").Value = 27
$ws.Range("G5").Value = 11
$ws.Range("L5").Value = 55
$ws.Range("M5").Value = 33
$ws.Range("R5").Value = 23
$ws.Range("S5").Value = 45
$ws.Range("T5").Value = 68

$ws = $wb.Worksheets.Item("great_50_exc")
$ws.Range("J2").Value = 59
$ws.Range("L2").Value = 32
$ws.Range("K3").Value = 75
$ws.Range("L3").Value = 42
$ws.Range("T3").Value = 77
$ws.Range("U3").Value = 57
$ws.Range("K4").Value = 55
$ws.Range("L4").Value = 34
$ws.Range("T4").Value = 48
$ws.Range("B5").Value = 26
$ws.Range("G5").Value = 11
$ws.Range("L5").Value = 31
$ws.Range("T5").Value = 37
